{"js": "// Update the worksheet date and the 3-digit-by-1-digit multiplication\n// problems. Replacements are applied in document order, one occurrence at a\n// time, so that a new value that happens to equal an old value used\n// elsewhere in the document (e.g. \"210\u00d77=\" is both an original cell and the\n// result of another cell's update) does not get re-matched by a later step.\nconst replacements = [\n  [\"2024-05-06 Monday\", \"2024-05-07 Tuesday\"],\n  [\"148\u00d77=\", \"378\u00d77=\"],\n  [\"210\u00d77=\", \"401\u00d77=\"],\n  [\"189\u00d75=\", \"468\u00d78=\"],\n  [\"887\u00d73=\", \"110\u00d76=\"],\n  [\"795\u00d73=\", \"414\u00d76=\"],\n  [\"983\u00d77=\", \"689\u00d79=\"],\n  [\"710\u00d78=\", \"662\u00d75=\"],\n  [\"292\u00d74=\", \"827\u00d75=\"],\n  [\"231\u00d72=\", \"972\u00d72=\"],\n  [\"735\u00d76=\", \"152\u00d78=\"],\n  [\"472\u00d79=\", \"661\u00d75=\"],\n  [\"495\u00d74=\", \"255\u00d72=\"],\n  [\"666\u00d78=\", \"547\u00d79=\"],\n  [\"374\u00d75=\", \"569\u00d79=\"],\n  [\"250\u00d76=\", \"653\u00d75=\"],\n  [\"874\u00d73=\", \"441\u00d73=\"],\n  [\"793\u00d74=\", \"533\u00d72=\"],\n  [\"821\u00d72=\", \"450\u00d72=\"],\n  [\"185\u00d78=\", \"225\u00d73=\"],\n  [\"135\u00d75=\", \"516\u00d78=\"],\n  [\"904\u00d72=\", \"534\u00d76=\"],\n  [\"283\u00d73=\", \"606\u00d77=\"],\n  [\"310\u00d76=\", \"210\u00d77=\"],\n  [\"246\u00d77=\", \"991\u00d77=\"],\n  [\"156\u00d75=\", \"219\u00d77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Only the first (document-order) occurrence is replaced; subsequent\n  // occurrences of the same old text (if any) are left for later steps.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 3-digit-by-1-digit multiplication\n# problems. Replacements are applied in document order, replacing only the\n# first (next) occurrence of each \"old\" value, so that a new value which\n# happens to equal an old value used elsewhere in the document (e.g.\n# \"210\u00d77=\" is both an original cell and the result of another cell's\n# update) does not get re-matched by a later step.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-05-06 Monday\", \"2024-05-07 Tuesday\"),\n    @(\"148\u00d77=\", \"378\u00d77=\"),\n    @(\"210\u00d77=\", \"401\u00d77=\"),\n    @(\"189\u00d75=\", \"468\u00d78=\"),\n    @(\"887\u00d73=\", \"110\u00d76=\"),\n    @(\"795\u00d73=\", \"414\u00d76=\"),\n    @(\"983\u00d77=\", \"689\u00d79=\"),\n    @(\"710\u00d78=\", \"662\u00d75=\"),\n    @(\"292\u00d74=\", \"827\u00d75=\"),\n    @(\"231\u00d72=\", \"972\u00d72=\"),\n    @(\"735\u00d76=\", \"152\u00d78=\"),\n    @(\"472\u00d79=\", \"661\u00d75=\"),\n    @(\"495\u00d74=\", \"255\u00d72=\"),\n    @(\"666\u00d78=\", \"547\u00d79=\"),\n    @(\"374\u00d75=\", \"569\u00d79=\"),\n    @(\"250\u00d76=\", \"653\u00d75=\"),\n    @(\"874\u00d73=\", \"441\u00d73=\"),\n    @(\"793\u00d74=\", \"533\u00d72=\"),\n    @(\"821\u00d72=\", \"450\u00d72=\"),\n    @(\"185\u00d78=\", \"225\u00d73=\"),\n    @(\"135\u00d75=\", \"516\u00d78=\"),\n    @(\"904\u00d72=\", \"534\u00d76=\"),\n    @(\"283\u00d73=\", \"606\u00d77=\"),\n    @(\"310\u00d76=\", \"210\u00d77=\"),\n    @(\"246\u00d77=\", \"991\u00d77=\"),\n    @(\"156\u00d75=\", \"219\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # wdFindContinue = 1, wdReplaceOne = 1 -- only the next occurrence is\n    # replaced, matching the sequential/document-order application below.\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n}\n"}
